$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matches source formatting)
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).NumberFormat = "@"

# Apply updated coin data
$ws.Cells.Item(2, 4).Value = '27.083.77'
$ws.Cells.Item(2, 5).Value = '  -0.45%  '

$ws.Cells.Item(3, 4).Value = '1.629.62'
$ws.Cells.Item(3, 5).Value = '  -1.22%  '

$ws.Cells.Item(4, 5).Value = '  -0.05%  '

$ws.Cells.Item(5, 4).Value = '216.10'
$ws.Cells.Item(5, 5).Value = '  -1.20%  '

$ws.Cells.Item(6, 4).Value = '0.513'
$ws.Cells.Item(6, 5).Value = '  +0.50%  '

$ws.Cells.Item(7, 5).Value = '  -0.06%  '

$ws.Cells.Item(8, 4).Value = '0.252'
$ws.Cells.Item(8, 5).Value = '  -1.54%  '

$ws.Cells.Item(9, 5).Value = '  -0.88%  '

$ws.Cells.Item(10, 5).Value = '  -0.93%  '

$ws.Cells.Item(11, 5).Value = '  -0.04%  '

$ws.Cells.Item(12, 4).Value = '1.857.42'
$ws.Cells.Item(12, 5).Value = '  -1.27%  '

$ws.Cells.Item(13, 4).Value = '1.627.68'
$ws.Cells.Item(13, 5).Value = '  -1.30%  '

$ws.Cells.Item(14, 4).Value = '4.10'
$ws.Cells.Item(14, 5).Value = '  -0.95%  '

$ws.Cells.Item(15, 4).Value = '0.540'
$ws.Cells.Item(15, 5).Value = '  +0.47%  '

$ws.Cells.Item(16, 4).Value = '65.72'
$ws.Cells.Item(16, 5).Value = '  -3.24%  '

$ws.Cells.Item(17, 4).Value = '27.056.19'
$ws.Cells.Item(17, 5).Value = '  -0.50%  '

$ws.Cells.Item(18, 5).Value = '  -0.69%  '

$ws.Cells.Item(19, 4).Value = '213.86'
$ws.Cells.Item(19, 5).Value = '  -3.02%  '

$ws.Cells.Item(20, 5).Value = '  -0.02%  '

$ws.Cells.Item(21, 5).Value = '  +1.16%  '

$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).Value = '4.38'
$ws.Cells.Item(22, 5).Value = '  -1.46%  '

$ws.Cells.Item(23, 2).Value = 'Toncoin'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(23, 4).Value = '2.50'
$ws.Cells.Item(23, 5).Value = '  +0.63%  '

$ws.Cells.Item(24, 5).Value = '  -1.90%  '

$ws.Cells.Item(25, 4).Value = '147.09'
$ws.Cells.Item(25, 5).Value = '  -0.83%  '

$ws.Cells.Item(26, 5).Value = '  -0.04%  '

$ws.Cells.Item(27, 4).Value = '7.37'
$ws.Cells.Item(27, 5).Value = '  -0.32%  '

$ws.Cells.Item(28, 5).Value = '  -1.50%  '

$ws.Cells.Item(29, 4).Value = '15.56'
$ws.Cells.Item(29, 5).Value = '  -1.84%  '

$ws.Cells.Item(30, 5).Value = '  -0.75%  '

$ws.Cells.Item(31, 5).Value = '  -1.04%  '

$ws.Cells.Item(32, 4).Value = '3.34'
$ws.Cells.Item(32, 5).Value = '  -0.50%  '

$ws.Cells.Item(33, 5).Value = '  -0.94%  '

$ws.Cells.Item(34, 4).Value = '1.299.42'
$ws.Cells.Item(34, 5).Value = '  +2.13%  '

$ws.Cells.Item(35, 5).Value = '  -1.36%  '

$ws.Cells.Item(36, 5).Value = '  -0.70%  '

$ws.Cells.Item(37, 5).Value = '  -1.49%  '

$ws.Cells.Item(38, 4).Value = '0.541'
$ws.Cells.Item(38, 5).Value = '  +0.10%  '

$ws.Cells.Item(39, 4).Value = '0.843'
$ws.Cells.Item(39, 5).Value = '  -0.12%  '

$ws.Cells.Item(40, 5).Value = '  -0.14%  '

$ws.Cells.Item(41, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(41, 4).Value = '0.807'
$ws.Cells.Item(41, 5).Value = '  -0.35%  '

$ws.Cells.Item(42, 2).Value = 'MXToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(42, 4).Value = '2.25'
$ws.Cells.Item(42, 5).Value = '  +4.16%  '

$ws.Cells.Item(43, 5).Value = '  -2.03%  '

$ws.Cells.Item(44, 4).Value = '1.766.72'
$ws.Cells.Item(44, 5).Value = '  -1.46%  '

$ws.Cells.Item(45, 4).Value = '61.99'

$ws.Cells.Item(46, 4).Value = '90.42'
$ws.Cells.Item(46, 5).Value = '  -2.16%  '

$ws.Cells.Item(47, 5).Value = '  -0.35%  '

$ws.Cells.Item(48, 4).Value = '0.0₆0105'
$ws.Cells.Item(48, 5).Value = '  -1.28%  '

$ws.Cells.Item(49, 4).Value = '0.0513'
$ws.Cells.Item(49, 5).Value = '  -0.36%  '

$ws.Cells.Item(50, 4).Value = '0.780'
$ws.Cells.Item(50, 5).Value = '  +16.39%  '

$ws.Cells.Item(51, 4).Value = '7.53'
$ws.Cells.Item(51, 5).Value = '  -2.61%  '
